$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I2").Value = 220.25
$ws.Range("K2").Value = 220.25
$ws.Range("M2").Value = -107.25
$ws.Range("H39").Value = 44.4
$ws.Range("I39").Value = 44.4
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 133.2
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = 162.8
$ws.Range("H42").Value = 8163.2856
$ws.Range("I42").Value = 1189.8334
$ws.Range("K42").Value = 3569.5002
$ws.Range("M42").Value = -3339.5002
$ws.Range("H70").Value = 100004530
$ws.Range("J70").Value = 125004910
$ws.Range("L70").Value = 375014730
$ws.Range("N70").Value = -375015270
$ws.Range("H73").Value = 100004530
$ws.Range("J73").Value = 125004910
$ws.Range("L73").Value = 375014730
$ws.Range("N73").Value = -375016602
$ws.Range("H80").Value = 35715612
$ws.Range("I80").Value = 83334056
$ws.Range("J80").Value = 1780.1875
$ws.Range("K80").Value = 250002168
$ws.Range("L80").Value = 5340.5625
$ws.Range("M80").Value = -250001170
$ws.Range("N80").Value = -7336.5625
$ws.Range("H83").Value = 35715612
$ws.Range("I83").Value = 83334056
$ws.Range("J83").Value = 1780.1875
$ws.Range("K83").Value = 750006504
$ws.Range("L83").Value = 16021.6875
$ws.Range("M83").Value = -750001512
$ws.Range("N83").Value = -26005.6875
$ws.Range("H92").Value = 1418.3
$ws.Range("I92").Value = 1669.3636
$ws.Range("J92").Value = 1111.4445
$ws.Range("K92").Value = 1669.3636
$ws.Range("L92").Value = 1111.4445
$ws.Range("M92").Value = -421.3635999999999
$ws.Range("N92").Value = -3607.4445
$ws.Range("H96").Value = 5707.2
$ws.Range("J96").Value = 11260
$ws.Range("L96").Value = 33780
$ws.Range("N96").Value = -36526
$ws.Range("H98").Value = 9958.032999999999
$ws.Range("I98").Value = 9872.286
$ws.Range("K98").Value = 9872.286
$ws.Range("M98").Value = -8374.286
$ws.Range("H113").Value = 3408.4849
$ws.Range("I113").Value = 3468.6538
$ws.Range("J113").Value = 3185
$ws.Range("K113").Value = 3468.6538
$ws.Range("L113").Value = 3185
$ws.Range("M113").Value = -214.6538
$ws.Range("N113").Value = -9693
$ws.Range("H118").Value = 1030.3846
$ws.Range("I118").Value = 932.25
$ws.Range("K118").Value = 2796.75
$ws.Range("M118").Value = -1139.75
$ws.Range("H122").Value = 9958.032999999999
$ws.Range("I122").Value = 9872.286
$ws.Range("K122").Value = 29616.858
$ws.Range("M122").Value = -27166.858
$ws.Range("H132").Value = 1924469
$ws.Range("I132").Value = 1152.425
$ws.Range("K132").Value = 3457.275
$ws.Range("M132").Value = -927.2749999999996
$ws.Range("H135").Value = 732.5333000000001
$ws.Range("I135").Value = 717.8148
$ws.Range("K135").Value = 6460.3332
$ws.Range("M135").Value = -3925.3332
$ws.Range("H137").Value = 532067.0600000001
$ws.Range("I137").Value = 1434.4166
$ws.Range("K137").Value = 4303.2498
$ws.Range("M137").Value = -1753.2498
$ws.Range("H138").Value = 2591404.5
$ws.Range("I138").Value = 1808
$ws.Range("J138").Value = 4843227.5
$ws.Range("K138").Value = 5424
$ws.Range("L138").Value = 14529682.5
$ws.Range("M138").Value = -284
$ws.Range("N138").Value = -14539962.5
$ws.Range("H141").Value = 5806.727
$ws.Range("I141").Value = 3517.9412
$ws.Range("K141").Value = 10553.8236
$ws.Range("M141").Value = -5373.8236
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5881.6
$ws.Range("I2").Value = 8870.833000000001
$ws.Range("J2").Value = 3888.7778
$ws.Range("K2").Value = 8870.833000000001
$ws.Range("L2").Value = 3888.7778
$ws.Range("M2").Value = -8757.833000000001
$ws.Range("N2").Value = -4114.7778
$ws.Range("H32").Value = 4934.164
$ws.Range("I32").Value = 4934.164
$ws.Range("K32").Value = 4934.164
$ws.Range("M32").Value = -4647.164
$ws.Range("H61").Value = 360563.75
$ws.Range("I61").Value = 3618.6562
$ws.Range("K61").Value = 3618.6562
$ws.Range("M61").Value = -3406.6562
$ws.Range("H62").Value = 420014050
$ws.Range("J62").Value = 420014050
$ws.Range("L62").Value = 420014050
$ws.Range("N62").Value = -420015298
$ws.Range("H63").Value = 10119.143
$ws.Range("I63").Value = 9617.546
$ws.Range("J63").Value = 10670.9
$ws.Range("K63").Value = 9617.546
$ws.Range("L63").Value = 10670.9
$ws.Range("M63").Value = -8931.546
$ws.Range("N63").Value = -12042.9
$ws.Range("H65").Value = 420014050
$ws.Range("J65").Value = 420014050
$ws.Range("L65").Value = 1260042150
$ws.Range("N65").Value = -1260048390
$ws.Range("H66").Value = 10119.143
$ws.Range("I66").Value = 9617.546
$ws.Range("J66").Value = 10670.9
$ws.Range("K66").Value = 48087.73
$ws.Range("L66").Value = 53354.5
$ws.Range("M66").Value = -44655.73
$ws.Range("N66").Value = -60218.5
$ws.Range("H74").Value = 1815.0238
$ws.Range("I74").Value = 1460.3939
$ws.Range("K74").Value = 1460.3939
$ws.Range("M74").Value = -586.3939
$ws.Range("H77").Value = 1815.0238
$ws.Range("I77").Value = 1460.3939
$ws.Range("K77").Value = 7301.9695
$ws.Range("M77").Value = -2933.9695
$ws.Range("H110").Value = 663.6
$ws.Range("I110").Value = 756.6667
$ws.Range("K110").Value = 756.6667
$ws.Range("M110").Value = 1288.3333
$ws.Range("H116").Value = 5881.6
$ws.Range("I116").Value = 8870.833000000001
$ws.Range("J116").Value = 3888.7778
$ws.Range("K116").Value = 8870.833000000001
$ws.Range("L116").Value = 3888.7778
$ws.Range("M116").Value = -6576.833000000001
$ws.Range("N116").Value = -8476.7778
$ws.Range("H136").Value = 360563.75
$ws.Range("I136").Value = 3618.6562
$ws.Range("K136").Value = 10855.9686
$ws.Range("M136").Value = -8305.9686
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 67998
$ws.Range("J2").Value = 74497
$ws.Range("L2").Value = 74497
$ws.Range("N2").Value = -74723
$ws.Range("H3").Value = 5881.6
$ws.Range("I3").Value = 8870.833000000001
$ws.Range("J3").Value = 3888.7778
$ws.Range("K3").Value = 8870.833000000001
$ws.Range("L3").Value = 3888.7778
$ws.Range("M3").Value = -8756.833000000001
$ws.Range("N3").Value = -4116.7778
$ws.Range("H29").Value = 858
$ws.Range("I29").Value = 858
$ws.Range("K29").Value = 858
$ws.Range("M29").Value = -569
$ws.Range("H43").Value = 261329.33
$ws.Range("J43").Value = 261329.33
$ws.Range("L43").Value = 261329.33
$ws.Range("N43").Value = -261691.33
$ws.Range("H86").Value = 55557010
$ws.Range("J86").Value = 166667420
$ws.Range("L86").Value = 166667420
$ws.Range("N86").Value = -166669666
$ws.Range("H89").Value = 55557010
$ws.Range("J89").Value = 166667420
$ws.Range("L89").Value = 833337100
$ws.Range("N89").Value = -833348332
$ws.Range("H94").Value = 1373.3429
$ws.Range("I94").Value = 1466.742
$ws.Range("K94").Value = 1466.742
$ws.Range("M94").Value = -1015.742
$ws.Range("H99").Value = 2589.5833
$ws.Range("I99").Value = 1619.4445
$ws.Range("K99").Value = 1619.4445
$ws.Range("M99").Value = -121.4445000000001
$ws.Range("H105").Value = 4352.154
$ws.Range("I105").Value = 3842
$ws.Range("K105").Value = 3842
$ws.Range("M105").Value = -2095
$ws.Range("H107").Value = 6849.7
$ws.Range("I107").Value = 5297.1816
$ws.Range("J107").Value = 8747.223
$ws.Range("K107").Value = 5297.1816
$ws.Range("L107").Value = 8747.223
$ws.Range("M107").Value = -3377.1816
$ws.Range("N107").Value = -12587.223
$ws.Range("H115").Value = 50000
$ws.Range("J115").Value = 50000
$ws.Range("L115").Value = 50000
$ws.Range("N115").Value = -53134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2719.9
$ws.Range("I2").Value = 2719.9
$ws.Range("K2").Value = 2719.9
$ws.Range("M2").Value = -2606.9
$ws.Range("H31").Value = 2899.2727
$ws.Range("I31").Value = 1677.762
$ws.Range("K31").Value = 1677.762
$ws.Range("M31").Value = -1382.762
$ws.Range("H32").Value = 425
$ws.Range("I32").Value = 425
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 425
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -109
$ws.Range("H34").Value = 2899.2727
$ws.Range("I34").Value = 1677.762
$ws.Range("K34").Value = 1677.762
$ws.Range("M34").Value = -1475.762
$ws.Range("H42").Value = 25987.334
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 25987.334
$ws.Range("K42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("M42").Value = 25987.334
$ws.Range("N42").Value = -27173.334
$ws.Range("H47").Value = 32250
$ws.Range("I47").Value = 30000
$ws.Range("K47").Value = 30000
$ws.Range("M47").Value = -29434
$ws.Range("H53").Value = 35000
$ws.Range("J53").Value = 35000
$ws.Range("L53").Value = 35000
$ws.Range("N53").Value = -36214
$ws.Range("H105").Value = 3770.389
$ws.Range("I105").Value = 4026.4167
$ws.Range("K105").Value = 4026.4167
$ws.Range("M105").Value = -2279.4167
$ws.Range("H122").Value = 13335828
$ws.Range("I122").Value = 2399.182
$ws.Range("K122").Value = 7197.545999999999
$ws.Range("M122").Value = -4747.545999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2299.8
$ws.Range("J34").Value = 1999.8
$ws.Range("L34").Value = 5999.4
$ws.Range("N34").Value = -6167.4
$ws.Range("H38").Value = 597.1579
$ws.Range("I38").Value = 28.5
$ws.Range("K38").Value = 85.5
$ws.Range("M38").Value = 261.5
$ws.Range("H82").Value = 49960.25
$ws.Range("I82").Value = 34566.668
$ws.Range("J82").Value = 59196.4
$ws.Range("K82").Value = 103700.004
$ws.Range("L82").Value = 177589.2
$ws.Range("M82").Value = -103294.004
$ws.Range("N82").Value = -178401.2
$ws.Range("H85").Value = 49960.25
$ws.Range("I85").Value = 34566.668
$ws.Range("J85").Value = 59196.4
$ws.Range("K85").Value = 103700.004
$ws.Range("L85").Value = 177589.2
$ws.Range("M85").Value = -102296.004
$ws.Range("N85").Value = -180397.2
$ws.Range("H97").Value = 166666860
$ws.Range("J97").Value = 269.66666
$ws.Range("L97").Value = 808.9999799999999
$ws.Range("N97").Value = -1800.99998
$ws.Range("H113").Value = 1330.1666
$ws.Range("J113").Value = 807.875
$ws.Range("L113").Value = 2423.625
$ws.Range("N113").Value = -6763.625
$ws.Range("H138").Value = 2568.1667
$ws.Range("I138").Value = 2665.8
$ws.Range("J138").Value = 2498.4285
$ws.Range("K138").Value = 7997.400000000001
$ws.Range("L138").Value = 7495.2855
$ws.Range("M138").Value = -2857.400000000001
$ws.Range("N138").Value = -17775.2855
$ws.Range("H141").Value = 76926110
$ws.Range("I141").Value = 90911624
$ws.Range("K141").Value = 272734872
$ws.Range("M141").Value = -272729692
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 32260944
$ws.Range("I80").Value = 62501650
$ws.Range("J80").Value = 4191.8
$ws.Range("K80").Value = 62501650
$ws.Range("L80").Value = 4191.8
$ws.Range("M80").Value = -62500652
$ws.Range("N80").Value = -6187.8
$ws.Range("H83").Value = 32260944
$ws.Range("I83").Value = 62501650
$ws.Range("J83").Value = 4191.8
$ws.Range("K83").Value = 312508250
$ws.Range("L83").Value = 20959
$ws.Range("M83").Value = -312503258
$ws.Range("N83").Value = -30943
$ws.Range("H102").Value = 1879.05
$ws.Range("I102").Value = 1202.7142
$ws.Range("J102").Value = 3457.1667
$ws.Range("K102").Value = 1202.7142
$ws.Range("L102").Value = 3457.1667
$ws.Range("M102").Value = 419.2858000000001
$ws.Range("N102").Value = -6701.1667
$ws.Range("H113").Value = 6234.1
$ws.Range("I113").Value = 1454
$ws.Range("J113").Value = 8808
$ws.Range("K113").Value = 1454
$ws.Range("L113").Value = 8808
$ws.Range("M113").Value = 716
$ws.Range("N113").Value = -13148
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 6666.3335
$ws.Range("J3").Value = 4999
$ws.Range("L3").Value = 4999
$ws.Range("N3").Value = -5223
$ws.Range("H6").Value = 27375
$ws.Range("J6").Value = 27375
$ws.Range("L6").Value = 27375
$ws.Range("N6").Value = -27599
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 1000
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = -828
$ws.Range("N14").Value = -1344
$ws.Range("H15").Value = 6666.3335
$ws.Range("J15").Value = 4999
$ws.Range("L15").Value = 4999
$ws.Range("N15").Value = -5339
$ws.Range("H17").Value = 4333
$ws.Range("I17").Value = 3250
$ws.Range("K17").Value = 3250
$ws.Range("M17").Value = -3080
$ws.Range("H20").Value = 504250
$ws.Range("I20").Value = 1000000
$ws.Range("J20").Value = 8500
$ws.Range("K20").Value = 1000000
$ws.Range("L20").Value = 8500
$ws.Range("M20").Value = -999774
$ws.Range("N20").Value = -8952
$ws.Range("H41").Value = 20000
$ws.Range("J41").Value = 20000
$ws.Range("L41").Value = 20000
$ws.Range("N41").Value = -20876
$ws.Range("H82").Value = 41669548
$ws.Range("I82").Value = 1159
$ws.Range("J82").Value = 62503744
$ws.Range("K82").Value = 1159
$ws.Range("L82").Value = 62503744
$ws.Range("M82").Value = -798
$ws.Range("N82").Value = -62504466
$ws.Range("H85").Value = 41669548
$ws.Range("I85").Value = 1159
$ws.Range("J85").Value = 62503744
$ws.Range("K85").Value = 1159
$ws.Range("L85").Value = 62503744
$ws.Range("M85").Value = 89
$ws.Range("N85").Value = -62506240
$ws.Range("H100").Value = 35718160
$ws.Range("I100").Value = 33335628
$ws.Range("J100").Value = 38467236
$ws.Range("K100").Value = 33335628
$ws.Range("L100").Value = 38467236
$ws.Range("M100").Value = -33335087
$ws.Range("N100").Value = -38468318
$ws.Range("H122").Value = 4512
$ws.Range("I122").Value = 3921.111
$ws.Range("J122").Value = 5176.75
$ws.Range("K122").Value = 11763.333
$ws.Range("L122").Value = 15530.25
$ws.Range("M122").Value = -9313.332999999999
$ws.Range("N122").Value = -20430.25
$ws.Range("H132").Value = 3273519.2
$ws.Range("I132").Value = 4633131.5
$ws.Range("K132").Value = 13899394.5
$ws.Range("M132").Value = -13896864.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 33496666
$ws.Range("J5").Value = 33496666
$ws.Range("L5").Value = 33496666
$ws.Range("N5").Value = -33496890
$ws.Range("H22").Value = 689.2857
$ws.Range("I22").Value = 570.8333
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = 570.8333
$ws.Range("L22").Value = 1400
$ws.Range("M22").Value = -277.8333
$ws.Range("N22").Value = -1986
$ws.Range("H41").Value = 18085.875
$ws.Range("J41").Value = 19318.2
$ws.Range("L41").Value = 19318.2
$ws.Range("N41").Value = -20098.2
$ws.Range("H55").Value = 13672.75
$ws.Range("J55").Value = 24897
$ws.Range("L55").Value = 24897
$ws.Range("N55").Value = -25451
$ws.Range("H100").Value = 742.86206
$ws.Range("I100").Value = 577.5
$ws.Range("K100").Value = 1155
$ws.Range("M100").Value = -614
$ws.Range("H113").Value = 1697.25
$ws.Range("I113").Value = 1866.3334
$ws.Range("J113").Value = 1190
$ws.Range("K113").Value = 5599.0002
$ws.Range("L113").Value = 3570
$ws.Range("M113").Value = -3429.0002
$ws.Range("N113").Value = -7910
$ws.Range("H132").Value = 15877084
$ws.Range("I132").Value = 19610692
$ws.Range("K132").Value = 58832076
$ws.Range("M132").Value = -58829546
$ws.Range("H141").Value = 94653.60000000001
$ws.Range("I141").Value = 94830
$ws.Range("J141").Value = 94609.5
$ws.Range("K141").Value = 94830
$ws.Range("L141").Value = 94609.5
$ws.Range("M141").Value = -89650
$ws.Range("N141").Value = -104969.5
